# Insert a new weekly record for "Arveja Verde" as row 36, pushing the
# existing rows 36:82 down to 37:83 (dimension grows from A1:R82 to A1:R83).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(36).Insert()

$ws.Cells.Item(36, 1).Value = 7
$ws.Cells.Item(36, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(36, 3).Value = "Ñuble"
$ws.Cells.Item(36, 4).Value = 44880
$ws.Cells.Item(36, 5).Value = 16
$ws.Cells.Item(36, 6).Value = 100112022
$ws.Cells.Item(36, 7).Value = "Arveja Verde"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 60
$ws.Cells.Item(36, 11).Value = 16000
$ws.Cells.Item(36, 12).Value = 17000
$ws.Cells.Item(36, 13).Value = 16500
$ws.Cells.Item(36, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(36, 15).Value = "Región del Maule"
$ws.Cells.Item(36, 16).Value = 660
$ws.Cells.Item(36, 17).Value = 25
$ws.Cells.Item(36, 18).Value = "Hortaliza"
